$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing unitary_measure (C) / total_measure (D) values so that
# each row satisfies: total_measure = units * unitary_measure
$ws.Range("C2").Value = 250
$ws.Range("D3").Value = 260
$ws.Range("D4").Value = 260
$ws.Range("D5").Value = 200
$ws.Range("D6").Value = 400
$ws.Range("C7").Value = 33.375
$ws.Range("D8").Value = 120
$ws.Range("C10").Value = 11.625
$ws.Range("C11").Value = 1.6
$ws.Range("C12").Value = 0.88
$ws.Range("D13").Value = 60
$ws.Range("C14").Value = 850
$ws.Range("C15").Value = 20
$ws.Range("D17").Value = 150

# Update the active cell selection to match the new state
$ws.Range("E9").Select()
